$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 (sldId 475) - "Content Placeholder 2" (shape id=3): update deadlines
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

$tr1.Replace(", deadline 19 May 2023", ", deadline 17 May 2024", 1, 0, 0) | Out-Null
$tr1.Replace("Final project, deadline 31 May 2023", "Final project, deadline 31 May 2024", 1, 0, 0) | Out-Null
$tr1.Replace("Exam, 17 June 2023", "Exam, 15 June 2024", 1, 0, 0) | Out-Null

# ---------------------------------------------------------------------------
# Slide 2 (sldId 703) - "Content Placeholder 2" (shape id=3):
#   * "ITX Flex must be enabled" -> "WISEflow Device Monitor must be enabled"
#     (split into 3 runs: "WISEflow" / " Device Monitor " / "must be enabled")
#   * "eksamen.au.dk" -> "wiseflow.au.dk"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

$para = $tr2.Paragraphs(3)

# "ITX Flex" is the first 8 characters of this paragraph (hyperlinked run) -
# turn it into " Device Monitor " while keeping its formatting/hyperlink.
$itxFlexRun = $para.Characters(1, 8)
$itxFlexRun.Text = " Device Monitor "

# Drop the now-duplicated leading space on the trailing "must be enabled" run.
$para.Replace(" must be enabled", "must be enabled", 1, 0, 0) | Out-Null

# Prepend a new "WISEflow" run (inherits the hyperlinked run's formatting)
# ahead of " Device Monitor ".
$insertionPoint = $para.Characters(1, 0)
$insertionPoint.InsertBefore("WISEflow") | Out-Null

# "eksamen.au.dk" -> "wiseflow.au.dk"
$tr2.Replace("eksamen.au.dk", "wiseflow.au.dk", 1, 0, 0) | Out-Null
